# Update countries & provincias Spain
# Refresh the "Pais" COVID-19 tracker sheet:
#  - bump the "Datos actualizados" timestamp
#  - refresh case counters for several countries
#  - Nepal overtakes Nigeria in the sorted list (same rows, swapped data)
#  - Bonaire jumps ahead of Santa Lucia / Timor Oriental / Nueva Caledonia

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: last-updated timestamp
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 15 de Septiembre de 2020 a las 14:12"

# Row 57: was Nigeria -> now Nepal (new, higher totals)
$ws.Cells.Item(57, 1).Value = "Nepal"
$ws.Cells.Item(57, 2).Value = 56788
$ws.Cells.Item(57, 3).Value = 1459
$ws.Cells.Item(57, 4).Value = 40638
$ws.Cells.Item(57, 5).Value = 15779
$ws.Cells.Item(57, 6).Value = 0
$ws.Cells.Item(57, 7).Value = 11
$ws.Cells.Item(57, 8).Value = 371

# Row 58: was Nepal -> now Nigeria (keeps Nigeria's previous totals)
$ws.Cells.Item(58, 1).Value = "Nigeria"
$ws.Cells.Item(58, 2).Value = 56388
$ws.Cells.Item(58, 3).Value = 0
$ws.Cells.Item(58, 4).Value = 44337
$ws.Cells.Item(58, 5).Value = 10968
$ws.Cells.Item(58, 6).Value = 0
$ws.Cells.Item(58, 7).Value = 0
$ws.Cells.Item(58, 8).Value = 1083

# Row 78: Bosnia y Herzegovina - refreshed totals
$ws.Cells.Item(78, 2).Value = 23929
$ws.Cells.Item(78, 3).Value = 294
$ws.Cells.Item(78, 4).Value = 16701
$ws.Cells.Item(78, 5).Value = 6503
$ws.Cells.Item(78, 7).Value = 20
$ws.Cells.Item(78, 8).Value = 725

# Row 81: Dinamarca - refreshed totals
$ws.Cells.Item(81, 2).Value = 20571
$ws.Cells.Item(81, 3).Value = 334
$ws.Cells.Item(81, 4).Value = 16557
$ws.Cells.Item(81, 5).Value = 3381

# Row 86: Madagascar - refreshed totals
$ws.Cells.Item(86, 2).Value = 15803
$ws.Cells.Item(86, 3).Value = 34
$ws.Cells.Item(86, 4).Value = 14452
$ws.Cells.Item(86, 5).Value = 1137
$ws.Cells.Item(86, 7).Value = 1
$ws.Cells.Item(86, 8).Value = 214

# Row 203: was Timor Oriental -> now Bonaire, San Eustaquio y Saba (new, higher totals)
$ws.Cells.Item(203, 1).Value = "Bonaire, San Eustaquio y Saba"
$ws.Cells.Item(203, 2).Value = 28
$ws.Cells.Item(203, 3).Value = 3
$ws.Cells.Item(203, 4).Value = 17
$ws.Cells.Item(203, 5).Value = 11

# Row 204: Santa Lucia stays put, values unchanged (kept for clarity)
$ws.Cells.Item(204, 1).Value = "Santa Lucia"
$ws.Cells.Item(204, 2).Value = 27
$ws.Cells.Item(204, 3).Value = 0
$ws.Cells.Item(204, 4).Value = 26
$ws.Cells.Item(204, 5).Value = 1

# Row 205: was Nueva Caledonia -> now Timor Oriental (keeps Timor's previous totals)
$ws.Cells.Item(205, 1).Value = "Timor Oriental"
$ws.Cells.Item(205, 2).Value = 27
$ws.Cells.Item(205, 3).Value = 0
$ws.Cells.Item(205, 4).Value = 26
$ws.Cells.Item(205, 5).Value = 1

# Row 206: was Bonaire, San Eustaquio y Saba -> now Nueva Caledonia (keeps its previous totals)
$ws.Cells.Item(206, 1).Value = "Nueva Caledonia"
$ws.Cells.Item(206, 2).Value = 26
$ws.Cells.Item(206, 3).Value = 0
$ws.Cells.Item(206, 4).Value = 26
$ws.Cells.Item(206, 5).Value = 0
